$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): add I0 and IF labels, matching the header style (s="1") used by H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-19
$i0 = @(9, 9, 8, 4, 7, 9, 7, 5, 6, 5, 7, 8, 9, 6, 8, 7, 3, 3)
$if = @(9, 9, 8, 5, 7, 9, 8, 6, 6, 5, 7, 8, 9, 7, 8, 7, 4, 4)

for ($r = 0; $r -lt $i0.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
